$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The old "notifications" sheet becomes "teacher_psswd": rename it in
#    place (keeps its sheetId/r:id) and replace its content with an empty
#    teacher-password table (mirrors the "student_pswd" sheet layout).
# ---------------------------------------------------------------------------
$teacherPsswd = $wb.Worksheets.Item("notifications")
$teacherPsswd.Name = "teacher_psswd"
$teacherPsswd.Cells.ClearContents()

$teacherPsswd.Range("A1").Value = "index number"
$teacherPsswd.Range("B1").Value = "password"
$teacherPsswd.Range("C1").Value = "First Name"

$teacherPsswd.Range("F3").Value = "num teachers"
$teacherPsswd.Range("F4").Value = 0

$teacherPsswd.Columns.AutoFit()

# ---------------------------------------------------------------------------
# 2. Add a brand-new sheet at the end of the workbook and name it
#    "notifications" - this is the new notifications table with updated
#    column headers and a sample enrollment notification row.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$notifications = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$notifications.Name = "notifications"

$notifications.Range("A1").Value = "Notification id"
$notifications.Range("B1").Value = "topic"
$notifications.Range("C1").Value = "sender"
$notifications.Range("D1").Value = "receiver"
$notifications.Range("E1").Value = "description"

$notifications.Range("A2").Value = 1
$notifications.Range("B2").Value = "Course Enrollment"
$notifications.Range("C2").Value = "Hasith"
$notifications.Range("D2").Value = "Sam Davis"
$notifications.Range("E2").Value = "Hasith would like to enroll in the Maths"

$notifications.Range("H2").Value = "num notifications"
$notifications.Range("H3").Value = 1

$notifications.Columns.AutoFit()

# ---------------------------------------------------------------------------
# 3. student_courses: bump the "No. of Students" counter (L4) from 1 to 2 -
#    a new student-course enrollment was recorded.
# ---------------------------------------------------------------------------
$studentCourses = $wb.Worksheets.Item("student_courses")
$studentCourses.Range("L4").Value = 2

# ---------------------------------------------------------------------------
# 4. View state: student_pswd now shows A1:C1 selected (its header row was
#    copied to seed the new teacher_psswd sheet); Teachers becomes the
#    active/selected tab with H4 selected.
# ---------------------------------------------------------------------------
$studentPswd = $wb.Worksheets.Item("student_pswd")
$studentPswd.Range("A1:C1").Select()

$teachers = $wb.Worksheets.Item("Teachers")
$teachers.Activate()
$teachers.Range("H4").Select()
